$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.150.88"
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = "'2.419.88"
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D5').Value = "'553.87"
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = "'0.590"
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = "'2.851.61"
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = "'60.064.57"
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = "'2.429.75"
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('D20').Value = "'327.70"
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = "'65.13"
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = "'0.178"
$ws.Range('E24').Value = '  +4.22%  '
$ws.Range('D25').Value = "'8.71"
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = "'1.40"
$ws.Range('E27').Value = '  +4.46%  '
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = "'170.28"
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = "'0.404"
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('D37').Value = "'4.22"
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').Value = "'326.88"
$ws.Range('E39').Value = '  +2.48%  '
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = "'144.55"
$ws.Range('E41').Value = '  +3.65%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = "'20.04"
$ws.Range('E43').Value = '  +2.71%  '
$ws.Range('D44').Value = "'0.0966"
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = "'11.05"
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('D51').Value = "'0.946"
$ws.Range('E51').Value = '  -0.52%  '
